$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Help Hours" labels in D10 and D11 to "Help Hours*"
$ws.Range("D10").Value = "Help Hours*"
$ws.Range("D11").Value = "Help Hours*"

# Reflect the updated active cell selection on the sheet
$ws.Activate()
$ws.Range("G16").Select()
